$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Snapshot existing values before the shift (old layout, columns A:G) ---
$old_A1 = $ws.Range("A1").Value2   # "i"
$old_B1 = $ws.Range("B1").Value2   # 10
$old_C1 = $ws.Range("C1").Value2   # 100
$old_D1 = $ws.Range("D1").Value2   # 1000
$old_E1 = $ws.Range("E1").Value2   # 10000
$old_G1 = $ws.Range("G1").Value2   # "runs: 1000"

$old_A2 = $ws.Range("A2").Value2   # "n"

$old_A3 = $ws.Range("A3").Value2   # 10
$old_B3 = $ws.Range("B3").Value2
$old_C3 = $ws.Range("C3").Value2
$old_D3 = $ws.Range("D3").Value2
$old_E3 = $ws.Range("E3").Value2

$old_A4 = $ws.Range("A4").Value2   # 100
$old_B4 = $ws.Range("B4").Value2

$old_A5 = $ws.Range("A5").Value2   # 1000
$old_B5 = $ws.Range("B5").Value2

$old_A6 = $ws.Range("A6").Value2   # 10000
$old_B6 = $ws.Range("B6").Value2

# --- Clear the old range so stale cells at e.g. G1/E6 don't linger ---
$ws.Range("A1:G6").ClearContents()

# --- Re-write the table shifted one column to the right (A:G -> B:H) ---
$ws.Range("B1").Value = $old_A1
$ws.Range("C1").Value = $old_B1
$ws.Range("D1").Value = $old_C1
$ws.Range("E1").Value = $old_D1
$ws.Range("F1").Value = $old_E1
$ws.Range("H1").Value = $old_G1

$ws.Range("B2").Value = $old_A2

$ws.Range("B3").Value = $old_A3
$ws.Range("C3").Value = $old_B3
$ws.Range("D3").Value = $old_C3
$ws.Range("E3").Value = $old_D3
$ws.Range("F3").Value = $old_E3

$ws.Range("B4").Value = $old_A4
$ws.Range("C4").Value = $old_B4

$ws.Range("B5").Value = $old_A5
$ws.Range("C5").Value = $old_B5

$ws.Range("B6").Value = $old_A6
$ws.Range("C6").Value = $old_B6

# --- New title for the first (Bucket sort) table ---
$ws.Range("A1").Value = "Bucket sort"

# Widen column A to fit the new titles (target stored width 14.140625 chars;
# the host quantizes column widths to 1/6-character steps, so 13.3 is the
# input that round-trips to the closest achievable stored width, 14.1667)
$ws.Columns.Item(1).ColumnWidth = 13.3

# --- New second table: Parallel merge (unfinished report) ---
$ws.Range("A8").Value = "Parellel merge"
$ws.Range("B8").Value = $old_A2
$ws.Range("C8").Value = "runtime"
$ws.Range("H8").Value = $old_G1

$ws.Range("B9").Value = 10
$ws.Range("B10").Value = 100
$ws.Range("B11").Value = 1000
$ws.Range("B12").Value = 10000

$ws.Range("G12").Select()
